$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.456.48'
$ws.Range("E2").Value = '  +3.71%  '
$ws.Range("D3").Value = '2.312.66'
$ws.Range("E3").Value = '  +3.03%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '105.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +9.82%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '309.02'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.16%  '
$ws.Range("E7").Value = '  +0.80%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +6.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.90'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.70'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0810'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("E13").Value = '  -1.04%  '
$ws.Range("E14").Value = '  +3.42%  '
$ws.Range("D15").Value = '2.672.14'
$ws.Range("E15").Value = '  +2.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.17'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.90%  '
$ws.Range("D17").Value = '2.316.11'
$ws.Range("E17").Value = '  +2.57%  '
$ws.Range("E18").Value = '  +3.10%  '
$ws.Range("D19").Value = '43.385.14'
$ws.Range("E19").Value = '  +3.72%  '
$ws.Range("E20").Value = '  -1.30%  '
$ws.Range("E21").Value = '  +2.74%  '
$ws.Range("E22").Value = '  +4.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.64'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.46%  '
$ws.Range("E25").Value = '  +4.57%  '
$ws.Range("E26").Value = '  +1.75%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.83'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.31'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.93%  '
$ws.Range("E30").Value = '  +1.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.10'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '162.27'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.00%  '
$ws.Range("E33").Value = '  +1.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.34'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.47%  '
$ws.Range("E37").Value = '  +2.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.58'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +14.07%  '
$ws.Range("E39").Value = '  -2.13%  '
$ws.Range("E40").Value = '  +3.41%  '
$ws.Range("E41").Value = '  +4.24%  '
$ws.Range("E42").Value = '  +0.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +14.90%  '
$ws.Range("E44").Value = '  +3.55%  '
$ws.Range("D45").Value = '1.966.46'
$ws.Range("E45").Value = '  +1.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.71'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("E47").Value = '  +7.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.28'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '58.06'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.93'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.58'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.91%  '
